$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view sizing (workbook.xml bookViews) ---
$excel.Left = -120
$excel.Top = -120
$excel.Width = 29040
$excel.Height = 15720

# --- Simple value / text edits ---
$ws.Range("C1").Value = "PCS0007"
$ws.Range("C2").Value = "Mr. Gaurab Pal"
$ws.Range("C3").Value = "L3"
$ws.Range("D4").Value = 41000
$ws.Range("E4").Value = 46000

# --- Formula edits ---
$ws.Range("D17").Formula = "=D14*12"
$ws.Range("E17").Formula = "=E14*12"

# --- Border/format swap between the "Insurance Benefits" header row (24)
#     and the "Leave" footer row (30) of the benefits box, via
#     copy/paste-special (format only) -- mirrors a Format Painter action.
#     Stash row 24's original formatting in a scratch range first so it
#     can be re-applied to row 30 after row 30's format is copied up.
$ws.Range("C24:D24").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("C30:D30").Copy()
$ws.Range("C24:D24").PasteSpecial(-4122)

$ws.Range("G1:H1").Copy()
$ws.Range("C30:D30").PasteSpecial(-4122)

$ws.Range("G1:H1").Clear()

# B1 picks up the same border treatment now applied to D24.
$ws.Range("D24").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# --- Selection change ---
$ws.Range("N14").Select() | Out-Null
